$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.308.18'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.006.31'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.97'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0768'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.17%  '
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.299.99'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.799'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.85'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.73%  '
$ws.Range('E16').Value = '  -4.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.000.73'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.162.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0834'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.85%  '
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.45%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('E29').Value = '  -6.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.119'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('E32').Value = '  -4.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0642'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.52'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.37'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.62%  '
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.36'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.33'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('E40').Value = '  +4.30%  '
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0929'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0212'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.408.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.46%  '
$ws.Range('E47').Value = '  -3.76%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.08%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.01'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.192.83'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('E51').Value = '  -7.48%  '
